$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "34.223.04"
Set-TextCell $ws.Range("E2") "  +1.13%  "

Set-TextCell $ws.Range("D3") "1.784.60"
Set-TextCell $ws.Range("E3") "  +0.31%  "

Set-TextCell $ws.Range("D5") "226.41"
Set-TextCell $ws.Range("E5") "  +0.80%  "

Set-TextCell $ws.Range("E6") "  +0.32%  "

Set-TextCell $ws.Range("E7") "  +0.10%  "

Set-TextCell $ws.Range("E8") "  -0.35%  "

Set-TextCell $ws.Range("E9") "  +0.80%  "

Set-TextCell $ws.Range("D10") "0.0694"
Set-TextCell $ws.Range("E10") "  +2.37%  "

Set-TextCell $ws.Range("D12") "2.042.08"
Set-TextCell $ws.Range("E12") "  +0.32%  "

Set-TextCell $ws.Range("B13") "Chainlink"
Set-TextCell $ws.Range("C13") "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell $ws.Range("D13") "11.08"
Set-TextCell $ws.Range("E13") "  -1.63%  "

Set-TextCell $ws.Range("B14") "WrappedEther"
Set-TextCell $ws.Range("C14") "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell $ws.Range("D14") "1.785.95"
Set-TextCell $ws.Range("E14") "  +0.05%  "

Set-TextCell $ws.Range("B15") "WrappedBTC"
Set-TextCell $ws.Range("C15") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextCell $ws.Range("D15") "34.192.46"
Set-TextCell $ws.Range("E15") "  +0.97%  "

Set-TextCell $ws.Range("B16") "Polygon"
Set-TextCell $ws.Range("C16") "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextCell $ws.Range("D16") "0.624"
Set-TextCell $ws.Range("E16") "  +2.10%  "

Set-TextCell $ws.Range("D17") "4.20"
Set-TextCell $ws.Range("E17") "  +1.26%  "

Set-TextCell $ws.Range("E18") "  +2.05%  "

Set-TextCell $ws.Range("D19") "0.0₃0808"
Set-TextCell $ws.Range("E19") "  +4.53%  "

Set-TextCell $ws.Range("D20") "247.20"
Set-TextCell $ws.Range("E20") "  +3.58%  "

Set-TextCell $ws.Range("D21") "10.99"

Set-TextCell $ws.Range("E22") "  +0.15%  "

Set-TextCell $ws.Range("E23") "  +2.20%  "

Set-TextCell $ws.Range("E24") "  -0.79%  "

Set-TextCell $ws.Range("D25") "162.76"
Set-TextCell $ws.Range("E25") "  +1.54%  "

Set-TextCell $ws.Range("E26") "  +2.37%  "

Set-TextCell $ws.Range("E27") "  +1.39%  "

Set-TextCell $ws.Range("E28") "  +1.44%  "

Set-TextCell $ws.Range("E30") "  +0.66%  "

Set-TextCell $ws.Range("E31") "  +1.87%  "

Set-TextCell $ws.Range("D32") "3.74"
Set-TextCell $ws.Range("E32") "  +4.39%  "

Set-TextCell $ws.Range("E33") "  +5.50%  "

Set-TextCell $ws.Range("E34") "  -1.38%  "

Set-TextCell $ws.Range("D35") "1.446.93"
Set-TextCell $ws.Range("E35") "  +4.39%  "

Set-TextCell $ws.Range("E36") "  +2.37%  "

Set-TextCell $ws.Range("D37") "2.42"
Set-TextCell $ws.Range("E37") "  +7.70%  "

Set-TextCell $ws.Range("E38") "  +3.55%  "

Set-TextCell $ws.Range("D39") "1.05"
Set-TextCell $ws.Range("E39") "  +0.72%  "

Set-TextCell $ws.Range("E40") "  +2.30%  "

Set-TextCell $ws.Range("E41") "  -1.26%  "

Set-TextCell $ws.Range("D42") "0.926"
Set-TextCell $ws.Range("E42") "  +1.70%  "

Set-TextCell $ws.Range("E43") "  +0.33%  "

Set-TextCell $ws.Range("E44") "  +1.02%  "

Set-TextCell $ws.Range("B45") "Kaspa"
Set-TextCell $ws.Range("C45") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws.Range("D45") "0.0510"
Set-TextCell $ws.Range("E45") "  +0.83%  "

Set-TextCell $ws.Range("B46") "FraxShare"
Set-TextCell $ws.Range("C46") "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws.Range("D46") "6.08"
Set-TextCell $ws.Range("E46") "  +3.89%  "

Set-TextCell $ws.Range("E47") "  -0.18%  "

Set-TextCell $ws.Range("D48") "0.0₆0135"
Set-TextCell $ws.Range("E48") "  -3.67%  "

Set-TextCell $ws.Range("D49") "1.943.78"
Set-TextCell $ws.Range("E49") "  +0.35%  "

Set-TextCell $ws.Range("E50") "  -2.21%  "

Set-TextCell $ws.Range("E51") "  +0.12%  "
